$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A24").Value = "2025-04-28 21:55:36"
$ws.Range("B24").Value = 64
